$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("H74").Value = 4657.643
$ws.Range("I74").Value = 4657.643
$ws.Range("K74").Value = 4657.643
$ws.Range("M74").Value = -3721.643
$ws.Range("H77").Value = 4657.643
$ws.Range("I77").Value = 4657.643
$ws.Range("K77").Value = 23288.215
$ws.Range("M77").Value = -18608.215
$ws.Range("H80").Value = 90919630
$ws.Range("I80").Value = 333334340
$ws.Range("J80").Value = 14123.5
$ws.Range("K80").Value = 1000003020
$ws.Range("L80").Value = 42370.5
$ws.Range("M80").Value = -1000002022
$ws.Range("N80").Value = -44366.5
$ws.Range("H83").Value = 90919630
$ws.Range("I83").Value = 333334340
$ws.Range("J83").Value = 14123.5
$ws.Range("K83").Value = 3000009060
$ws.Range("L83").Value = 127111.5
$ws.Range("M83").Value = -3000004068
$ws.Range("N83").Value = -137095.5
$ws.Range("H98").Value = 1336.5106
$ws.Range("I98").Value = 1367.25
$ws.Range("J98").Value = 885.6667
$ws.Range("K98").Value = 1367.25
$ws.Range("L98").Value = 885.6667
$ws.Range("M98").Value = 130.75
$ws.Range("N98").Value = -3881.6667
$ws.Range("H116").Value = 21799.4
$ws.Range("I116").Value = 30165.666
$ws.Range("K116").Value = 30165.666
$ws.Range("M116").Value = -26723.666
$ws.Range("H122").Value = 1336.5106
$ws.Range("I122").Value = 1367.25
$ws.Range("J122").Value = 885.6667
$ws.Range("K122").Value = 4101.75
$ws.Range("L122").Value = 2657.0001
$ws.Range("M122").Value = -1651.75
$ws.Range("N122").Value = -7557.0001
$ws.Range("H123").Value = 70747.25
$ws.Range("J123").Value = 70747.25
$ws.Range("L123").Value = 70747.25
$ws.Range("N123").Value = -80547.25
$ws.Range("H129").Value = 1982.75
$ws.Range("J129").Value = 2472.7273
$ws.Range("L129").Value = 7418.1819
$ws.Range("N129").Value = -17418.1819
$ws.Range("H132").Value = 5599.355
$ws.Range("I132").Value = 6221.125
$ws.Range("K132").Value = 18663.375
$ws.Range("M132").Value = -16133.375
$ws.Range("H133").Value = 124997.5
$ws.Range("J133").Value = 124997.5
$ws.Range("L133").Value = 124997.5
$ws.Range("N133").Value = -135117.5
$ws.Range("H137").Value = 71366
$ws.Range("I137").Value = 104124.5
$ws.Range("J137").Value = 5849
$ws.Range("K137").Value = 312373.5
$ws.Range("L137").Value = 17547
$ws.Range("M137").Value = -309823.5
$ws.Range("N137").Value = -22647
$ws.Range("H138").Value = 3084.45
$ws.Range("J138").Value = 3422.7317
$ws.Range("L138").Value = 10268.1951
$ws.Range("N138").Value = -20548.1951
$ws.Range("H141").Value = 3641.7144
$ws.Range("I141").Value = 1939
$ws.Range("K141").Value = 5817
$ws.Range("M141").Value = -637

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 4601.3335
$ws.Range("I36").Value = 4601.3335
$ws.Range("K36").Value = 4601.3335
$ws.Range("M36").Value = -4255.3335
$ws.Range("H45").Value = 4115.1875
$ws.Range("J45").Value = 4892
$ws.Range("L45").Value = 4892
$ws.Range("N45").Value = -5646
$ws.Range("H61").Value = 3103.238
$ws.Range("I61").Value = 2830.3333
$ws.Range("J61").Value = 4740.6665
$ws.Range("K61").Value = 2830.3333
$ws.Range("L61").Value = 4740.6665
$ws.Range("M61").Value = -2618.3333
$ws.Range("N61").Value = -5164.6665
$ws.Range("H74").Value = 2259.6155
$ws.Range("I74").Value = 2259.6155
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 2259.6155
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -1385.6155
$ws.Range("H77").Value = 2259.6155
$ws.Range("I77").Value = 2259.6155
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 11298.0775
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -6930.077499999999
$ws.Range("H102").Value = 1513.0714
$ws.Range("I102").Value = 1358.9
$ws.Range("J102").Value = 1898.5
$ws.Range("K102").Value = 1358.9
$ws.Range("L102").Value = 1898.5
$ws.Range("M102").Value = 263.0999999999999
$ws.Range("N102").Value = -5142.5
$ws.Range("H105").Value = 131499.5
$ws.Range("J105").Value = 131499.5
$ws.Range("L105").Value = 131499.5
$ws.Range("N105").Value = -138487.5
$ws.Range("H136").Value = 3103.238
$ws.Range("I136").Value = 2830.3333
$ws.Range("J136").Value = 4740.6665
$ws.Range("K136").Value = 8490.999899999999
$ws.Range("L136").Value = 14221.9995
$ws.Range("M136").Value = -5940.999899999999
$ws.Range("N136").Value = -19321.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H60").Value = 80104.25
$ws.Range("J60").Value = 80104.25
$ws.Range("L60").Value = 80104.25
$ws.Range("N60").Value = -81302.25
$ws.Range("H88").Value = 100000
$ws.Range("J88").Value = 100000
$ws.Range("L88").Value = 100000
$ws.Range("N88").Value = -100812
$ws.Range("H91").Value = 100000
$ws.Range("J91").Value = 100000
$ws.Range("L91").Value = 100000
$ws.Range("N91").Value = -102808
$ws.Range("H94").Value = 726.6539
$ws.Range("I94").Value = 643.3333
$ws.Range("J94").Value = 914.125
$ws.Range("K94").Value = 643.3333
$ws.Range("L94").Value = 914.125
$ws.Range("M94").Value = -192.3333
$ws.Range("N94").Value = -1816.125
$ws.Range("H99").Value = 2553
$ws.Range("I99").Value = 1342.75
$ws.Range("K99").Value = 1342.75
$ws.Range("M99").Value = 155.25
$ws.Range("H105").Value = 3666.4348
$ws.Range("I105").Value = 2297.4285
$ws.Range("K105").Value = 2297.4285
$ws.Range("M105").Value = -550.4285
$ws.Range("H107").Value = 2008.9333
$ws.Range("I107").Value = 2081.0715
$ws.Range("K107").Value = 2081.0715
$ws.Range("M107").Value = -161.0715
$ws.Range("H134").Value = 11114225
$ws.Range("I134").Value = 1255.381
$ws.Range("K134").Value = 3766.143
$ws.Range("M134").Value = -1231.143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3279.6333
$ws.Range("I31").Value = 2105.05
$ws.Range("J31").Value = 5628.8
$ws.Range("K31").Value = 2105.05
$ws.Range("L31").Value = 5628.8
$ws.Range("M31").Value = -1810.05
$ws.Range("N31").Value = -6218.8
$ws.Range("H34").Value = 3279.6333
$ws.Range("I34").Value = 2105.05
$ws.Range("J34").Value = 5628.8
$ws.Range("K34").Value = 2105.05
$ws.Range("L34").Value = 5628.8
$ws.Range("M34").Value = -1903.05
$ws.Range("N34").Value = -6032.8
$ws.Range("H60").Value = 46844.668
$ws.Range("J60").Value = 46844.668
$ws.Range("L60").Value = 46844.668
$ws.Range("N60").Value = -47866.668
$ws.Range("H105").Value = 1512.2667
$ws.Range("I105").Value = 1421.9231
$ws.Range("J105").Value = 2099.5
$ws.Range("K105").Value = 1421.9231
$ws.Range("L105").Value = 2099.5
$ws.Range("M105").Value = 325.0769
$ws.Range("N105").Value = -5593.5
$ws.Range("H132").Value = 1680.5264
$ws.Range("I132").Value = 1196
$ws.Range("J132").Value = 2511.1428
$ws.Range("K132").Value = 3588
$ws.Range("L132").Value = 7533.428400000001
$ws.Range("M132").Value = -1058
$ws.Range("N132").Value = -12593.4284
$ws.Range("H134").Value = 1707.0444
$ws.Range("I134").Value = 1143.4857
$ws.Range("K134").Value = 3430.4571
$ws.Range("M134").Value = -895.4570999999996

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 114688670
$ws.Range("J4").Value = 177553230
$ws.Range("L4").Value = 532659690
$ws.Range("N4").Value = -532659914
$ws.Range("H132").Value = 1432.909
$ws.Range("I132").Value = 794.9231
$ws.Range("J132").Value = 2354.4443
$ws.Range("K132").Value = 7154.3079
$ws.Range("L132").Value = 21189.9987
$ws.Range("M132").Value = -4624.3079
$ws.Range("N132").Value = -26249.9987

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2910.4736
$ws.Range("I102").Value = 2023.0769
$ws.Range("J102").Value = 4833.1665
$ws.Range("K102").Value = 2023.0769
$ws.Range("L102").Value = 4833.1665
$ws.Range("M102").Value = -401.0769
$ws.Range("N102").Value = -8077.1665
$ws.Range("H111").Value = 60992
$ws.Range("J111").Value = 60992
$ws.Range("L111").Value = 60992
$ws.Range("N111").Value = -67126

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3295
$ws.Range("I46").Value = 665
$ws.Range("K46").Value = 665
$ws.Range("M46").Value = -477
$ws.Range("H122").Value = 2984.0833
$ws.Range("I122").Value = 2455.4443
$ws.Range("K122").Value = 7366.3329
$ws.Range("M122").Value = -4916.3329
$ws.Range("H136").Value = 3654.879
$ws.Range("I136").Value = 3495.5908
$ws.Range("J136").Value = 3973.4546
$ws.Range("K136").Value = 10486.7724
$ws.Range("L136").Value = 11920.3638
$ws.Range("M136").Value = -7936.7724
$ws.Range("N136").Value = -17020.3638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 21331
$ws.Range("I32").Value = 11285.2
$ws.Range("K32").Value = 11285.2
$ws.Range("M32").Value = -10968.2
$ws.Range("H81").Value = 3682.5217
$ws.Range("I81").Value = 2468.9
$ws.Range("J81").Value = 4616.077
$ws.Range("K81").Value = 4937.8
$ws.Range("L81").Value = 9232.154
$ws.Range("M81").Value = -3876.8
$ws.Range("N81").Value = -11354.154
$ws.Range("H84").Value = 3682.5217
$ws.Range("I84").Value = 2468.9
$ws.Range("J84").Value = 4616.077
$ws.Range("K84").Value = 24689
$ws.Range("L84").Value = 46160.77
$ws.Range("M84").Value = -19385
$ws.Range("N84").Value = -56768.77
$ws.Range("H132").Value = 3692.3635
$ws.Range("I132").Value = 3354.24
$ws.Range("K132").Value = 10062.72
$ws.Range("M132").Value = -7532.719999999999

# Remove cells that should no longer exist
$wb.Worksheets.Item("ALC").Range("M10").ClearContents()
$wb.Worksheets.Item("ARM").Range("N74").ClearContents()
$wb.Worksheets.Item("ARM").Range("N77").ClearContents()

Write-Output "Done applying updates"